$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Scenarios")

# Fill in the new "ARE_GUI_7" stresstest row (row 13), which was previously empty.
$ws.Range("A13").Value = "ARE_GUI_7"
$ws.Range("B13").Value = "Stresstest Start model"
$ws.Range("C13").Value = "Test ARE_GUI_2"
$ws.Range("D13").Value = "1. Execute Test ARE_GUI_2 by clicking 10 times onto 'Start' button as fast as possible"
$ws.Range("E13").Value = "The model must be started 10 times sequentially and successfully`nThe last model start must have a clean state and must not have orphaned GUI elements in the ARE GUI  panel`nThe ARE must not crash"

# Update the active cell selection to match the final state recorded in the file.
$ws.Range("E14").Select()
